{"js": "// Update the date line and every division-problem cell in the practice\n// table. The document is a single title paragraph followed by a 5x20\n// table; only every 4th row actually contains a \"a\u00f7b=\" problem text, the\n// remaining rows are blank answer rows. We walk body.paragraphs in\n// document order (title paragraph first, then every table-cell paragraph\n// row by row, column by column) and overwrite just the non-blank ones\n// positionally, so the edit is robust even though several of the new\n// values collide with old values used elsewhere in the document.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Ordered replacements, matching the document's paragraph order exactly\n// (title, then each table cell in row-major order, skipping blank rows).\nconst replacements = [\n  \"2024-08-25 Sunday\",\n  \"32\u00f75=\",\n  \"31\u00f75=\",\n  \"58\u00f75=\",\n  \"23\u00f72=\",\n  \"43\u00f74=\",\n  \"56\u00f76=\",\n  \"18\u00f74=\",\n  \"32\u00f76=\",\n  \"43\u00f74=\",\n  \"35\u00f77=\",\n  \"78\u00f75=\",\n  \"63\u00f74=\",\n  \"49\u00f74=\",\n  \"50\u00f77=\",\n  \"33\u00f72=\",\n  \"93\u00f72=\",\n  \"27\u00f72=\",\n  \"58\u00f79=\",\n  \"74\u00f74=\",\n  \"38\u00f72=\",\n  \"89\u00f72=\",\n  \"96\u00f78=\",\n  \"47\u00f73=\",\n  \"59\u00f77=\",\n  \"45\u00f76=\"\n];\n\nlet r = 0;\nfor (let i = 0; i < paragraphs.items.length && r < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text !== \"\") {\n    para.insertText(replacements[r], Word.InsertLocation.replace);\n    r++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every division-problem cell in the practice\n# table. The document is a single title paragraph followed by a 5x20\n# table; only every 4th row actually contains a \"a\u00f7b=\" problem text, the\n# remaining rows are blank answer rows. We walk $d.Paragraphs in document\n# order (title paragraph first, then every table-cell paragraph row by\n# row, column by column) and overwrite just the non-blank ones\n# positionally, so the edit is robust even though several of the new\n# values collide with old values used elsewhere in the document.\n\n$d = $word.ActiveDocument\n\n# Ordered replacements, matching the document's paragraph order exactly\n# (title, then each table cell in row-major order, skipping blank rows).\n$replacements = @(\n    \"2024-08-25 Sunday\",\n    \"32\u00f75=\",\n    \"31\u00f75=\",\n    \"58\u00f75=\",\n    \"23\u00f72=\",\n    \"43\u00f74=\",\n    \"56\u00f76=\",\n    \"18\u00f74=\",\n    \"32\u00f76=\",\n    \"43\u00f74=\",\n    \"35\u00f77=\",\n    \"78\u00f75=\",\n    \"63\u00f74=\",\n    \"49\u00f74=\",\n    \"50\u00f77=\",\n    \"33\u00f72=\",\n    \"93\u00f72=\",\n    \"27\u00f72=\",\n    \"58\u00f79=\",\n    \"74\u00f74=\",\n    \"38\u00f72=\",\n    \"89\u00f72=\",\n    \"96\u00f78=\",\n    \"47\u00f73=\",\n    \"59\u00f77=\",\n    \"45\u00f76=\"\n)\n\n$r = 0\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $trimmed = $t.TrimEnd([char]13, [char]7)\n    if ($trimmed.Length -gt 0 -and $r -lt $replacements.Length) {\n        $p.Range.Text = $replacements[$r]\n        $r++\n    }\n}\n"}
